$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-18 (columns A, B, C)
$data = @(
    @("<your>",   "<your>",   21),
    @("<kilo>",   "<kilo>",   17),
    @("<down>",   "<down>",   22),
    @("<quebec>", "<quebec>", 20),
    @("<down>",   "<down>",   22),
    @("<other>",  "<other>",  26),
    @("<four>",   "<for>",    19),
    @("<oh>",     "<oh>",     29),
    @("<all>",    "<all>",    21),
    @("<been>",   "<been>",   17),
    @("<can>",    "<come>",   19),
    @("<some>",   "<some>",   21),
    @("<see>",    "<see>",    19),
    @("<it>",     "<it>",     20),
    @("<paste>",  "<paste>",  22),
    @("<zero>",   "<zero>",   29),
    @("<for>",    "<for>",    25)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
